$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 264
$ws.Range("I2").Value = 737
$ws.Range("J2").Value = 2817
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 822
$ws.Range("M2").Value = 43
$ws.Range("N2").Value = 539
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 8
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 36
$ws.Range("S2").Value = 284
$ws.Range("T2").Value = 505
$ws.Range("U2").Value = 36
$ws.Range("V2").Value = 4523
$ws.Range("X2").Value = 4328
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 29
